# Re-sort the "EC" detail table (rows 16-21) on sheet "Hoja1": the
# "Periodo Mora" column (E) is reversed from ascending (2404..2409) to
# descending (2409..2404), and the paired "Valor Mora" column (F)
# travels together with each period, so the previously-mismatched value
# (32933 on the last period) now sits with period 2409 in the first row,
# and the common value (52000) ends up on period 2404 in the last row.
#
# Before:                              After:
#   E16=2404 F16=52000                   E16=2409 F16=32933
#   E17=2405 F17=52000                   E17=2408 F17=52000
#   E18=2406 F18=52000                   E18=2407 F18=52000
#   E19=2407 F19=52000                   E19=2406 F19=52000
#   E20=2408 F20=52000                   E20=2405 F20=52000
#   E21=2409 F21=32933                   E21=2404 F21=52000

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$periodos = @("2409", "2408", "2407", "2406", "2405", "2404")
$valores  = @(32933, 52000, 52000, 52000, 52000, 52000)

for ($i = 0; $i -lt 6; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodos[$i]
    $ws.Range("F$row").Value = $valores[$i]
}
